# Updated symbol list on Sat Dec 17 08:50:15 UTC 2022 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text cells (coin names, links, volume labels) - safe to assign directly,
# Excel will not reinterpret these as numbers.
$textCells = @{
    "B10" = "One"
    "C10" = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
    "E10" = "9OneONE"

    "B11" = "WazirX"
    "C11" = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
    "E11" = "10WazirXWRX"

    "B12" = "MandalaExchangeToken"
    "C12" = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
    "E12" = "11MandalaExchangeTokenMDX"

    "B13" = "LiechtensteinCryptoassetsExchange"
    "C13" = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
    "E13" = "12LiechtensteinCryptoassetsExchangeLCX"

    "B14" = "BitrueCoin"
    "C14" = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
    "E14" = "13BitrueCoinBTR"

    "B15" = "BitMartToken"
    "C15" = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
    "E15" = "14BitMartTokenBMX"

    "B16" = "BitForexToken"
    "C16" = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
    "E16" = "15BitForexTokenBF"

    "B17" = "MCDex"
    "C17" = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
    "E17" = "16MCDexMCB"

    "B18" = "CoinExToken"
    "C18" = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
    "E18" = "17CoinExTokenCET"

    "E23" = "22LEOLEO"

    "E41" = "40KickTokenKICKBestin24h"

    "B42" = "CEJI"
    "C42" = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
    "E42" = "41CEJICEJI"

    "B43" = "BKEXToken"
    "C43" = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
    "E43" = "42BKEXTokenBKK"
}

foreach ($key in $textCells.Keys) {
    $ws.Range($key).Value = $textCells[$key]
}

# Price cells (column D) hold numeric-looking strings but must remain stored as
# literal text (matching the source's inline-string/general-format cells, and
# preserving exact formatting such as trailing zeros, e.g. "0.07408"). Forcing
# the NumberFormat to Text ("@") before assignment keeps the plain-text value
# instead of Excel parsing it into a float; resetting the Style back to
# "Normal" afterwards avoids leaving a stray number-format style applied to
# the cell.
$priceCells = @{
    "D2"  = "235.62"
    "D3"  = "22.32"
    "D4"  = "5.414"
    "D5"  = "0.05635"
    "D6"  = "3.374"
    "D7"  = "6.480"
    "D8"  = "1.073"
    "D9"  = "0.7826"
    "D10" = "0.0005732"
    "D11" = "0.1398"
    "D12" = "0.07408"
    "D13" = "0.03191"
    "D14" = "0.02956"
    "D15" = "0.09261"
    "D16" = "0.001662"
    "D17" = "3.264"
    "D18" = "0.04760"
    "D19" = "0.006212"
    "D20" = "0.005110"
    "D23" = "3.895"
    "D24" = "2.147"
    "D27" = "0.0004992"
    "D40" = "0.04047"
    "D41" = "0.007027"
    "D42" = "0.003502"
    "D43" = "0.1039"
    "D44" = "0.009293"
    "D45" = "0.00005442"
    "D47" = "0.6755"
    "D48" = "0.03963"
    "D49" = "0.00002101"
}

foreach ($key in $priceCells.Keys) {
    $ws.Range($key).NumberFormat = "@"
    $ws.Range($key).Value = $priceCells[$key]
    $ws.Range($key).Style = "Normal"
}
